# Auto-generated edit script applying the Phantom_Profits.xlsx cell-value diff
# (scheduled runner refresh of market-price-derived Leve profit columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 858.2368
$ws.Range("I15").Value = 858.2368
$ws.Range("K15").Value = 2574.7104
$ws.Range("M15").Value = -2405.7104
$ws.Range("H63").Value = 45000
$ws.Range("I63").Value = 45000
$ws.Range("K63").Value = 45000
$ws.Range("M63").Value = -44376
$ws.Range("H66").Value = 45000
$ws.Range("I66").Value = 45000
$ws.Range("K66").Value = 135000
$ws.Range("M66").Value = -131880
$ws.Range("H137").Value = 9805332
$ws.Range("I137").Value = 15152445
$ws.Range("K137").Value = 45457335
$ws.Range("M137").Value = -45454785
$ws.Range("H141").Value = 2299.8667
$ws.Range("I141").Value = 1958
$ws.Range("J141").Value = 3667.3333
$ws.Range("K141").Value = 5874
$ws.Range("L141").Value = 11001.9999
$ws.Range("M141").Value = -694
$ws.Range("N141").Value = -21361.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 2420.7
$ws.Range("I4").Value = 3116.1667
$ws.Range("J4").Value = 1377.5
$ws.Range("K4").Value = 3116.1667
$ws.Range("L4").Value = 1377.5
$ws.Range("M4").Value = -3000.1667
$ws.Range("N4").Value = -1609.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 17547.555
$ws.Range("I82").Value = 17547.555
$ws.Range("K82").Value = 17547.555
$ws.Range("M82").Value = -17164.555
$ws.Range("H85").Value = 17547.555
$ws.Range("I85").Value = 17547.555
$ws.Range("K85").Value = 17547.555
$ws.Range("M85").Value = -16221.555
$ws.Range("H86").Value = 26505.725
$ws.Range("I86").Value = 25454.79
$ws.Range("K86").Value = 25454.79
$ws.Range("M86").Value = -24331.79
$ws.Range("H89").Value = 26505.725
$ws.Range("I89").Value = 25454.79
$ws.Range("K89").Value = 127273.95
$ws.Range("M89").Value = -121657.95
$ws.Range("H97").Value = 9484.666999999999
$ws.Range("I97").Value = 9484.666999999999
$ws.Range("K97").Value = 9484.666999999999
$ws.Range("M97").Value = -8493.666999999999
$ws.Range("H107").Value = 3686.5
$ws.Range("I107").Value = 3782.7144
$ws.Range("J107").Value = 3013
$ws.Range("K107").Value = 3782.7144
$ws.Range("L107").Value = 3013
$ws.Range("M107").Value = -1862.7144
$ws.Range("N107").Value = -6853

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 329.5
$ws.Range("I5").Value = 339.5
$ws.Range("J5").Value = 299.5
$ws.Range("K5").Value = 339.5
$ws.Range("L5").Value = 299.5
$ws.Range("M5").Value = -227.5
$ws.Range("N5").Value = -523.5
$ws.Range("H22").Value = 337.6842
$ws.Range("I22").Value = 360.46155
$ws.Range("J22").Value = 288.33334
$ws.Range("K22").Value = 360.46155
$ws.Range("L22").Value = 288.33334
$ws.Range("M22").Value = -10.46154999999999
$ws.Range("N22").Value = -988.33334
$ws.Range("H111").Value = 84495
$ws.Range("J111").Value = 84495
$ws.Range("L111").Value = 84495
$ws.Range("N111").Value = -92675
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H141").Value = 349350.88
$ws.Range("J141").Value = 386769.75
$ws.Range("L141").Value = 386769.75
$ws.Range("N141").Value = -397129.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 3965.5557
$ws.Range("J9").Value = 3961.25
$ws.Range("L9").Value = 11883.75
$ws.Range("N9").Value = -12331.75
$ws.Range("H11").Value = 1073.2858
$ws.Range("I11").Value = 1235
$ws.Range("K11").Value = 3705
$ws.Range("M11").Value = -3565
$ws.Range("H68").Value = 977
$ws.Range("I68").Value = 977
$ws.Range("K68").Value = 2931
$ws.Range("M68").Value = -2120
$ws.Range("H71").Value = 977
$ws.Range("I71").Value = 977
$ws.Range("K71").Value = 8793
$ws.Range("M71").Value = -4737
$ws.Range("H131").Value = 1778.25
$ws.Range("I131").Value = 1379.8334
$ws.Range("J131").Value = 2973.5
$ws.Range("K131").Value = 4139.5002
$ws.Range("L131").Value = 8920.5
$ws.Range("M131").Value = 900.4997999999996
$ws.Range("N131").Value = -19000.5
$ws.Range("H140").Value = 529404.1
$ws.Range("I140").Value = 529404.1
$ws.Range("K140").Value = 1588212.3
$ws.Range("M140").Value = -1583032.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H44").Value = 28514
$ws.Range("I44").Value = 21997.5
$ws.Range("K44").Value = 21997.5
$ws.Range("M44").Value = -21401.5
$ws.Range("H101").Value = 34194.25
$ws.Range("J101").Value = 36507.715
$ws.Range("L101").Value = 36507.715
$ws.Range("N101").Value = -42997.715
$ws.Range("H126").Value = 7338.1
$ws.Range("I126").Value = 7300.067
$ws.Range("J126").Value = 7452.2
$ws.Range("K126").Value = 21900.201
$ws.Range("L126").Value = 22356.6
$ws.Range("M126").Value = -19430.201
$ws.Range("N126").Value = -27296.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3896.6
$ws.Range("I46").Value = 2061.8
$ws.Range("K46").Value = 2061.8
$ws.Range("M46").Value = -1873.8
$ws.Range("H127").Value = 100000
$ws.Range("J127").Value = 100000
$ws.Range("L127").Value = 100000
$ws.Range("N127").Value = -109920

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 17000
$ws.Range("I41").Value = 17000
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 17000
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -16610
$ws.Range("N41").ClearContents()
$ws.Range("H52").Value = 19299.666
$ws.Range("I52").Value = 19299.666
$ws.Range("K52").Value = 19299.666
$ws.Range("M52").Value = -19073.666
$ws.Range("H54").Value = 35896.9
$ws.Range("I54").Value = 9394
$ws.Range("J54").Value = 62399.8
$ws.Range("K54").Value = 9394
$ws.Range("L54").Value = 62399.8
$ws.Range("M54").Value = -8874
$ws.Range("N54").Value = -63439.8
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("N56").ClearContents()
$ws.Range("H62").Value = 6250
$ws.Range("I62").Value = 5500
$ws.Range("K62").Value = 5500
$ws.Range("M62").Value = -4876
$ws.Range("H65").Value = 6250
$ws.Range("I65").Value = 5500
$ws.Range("K65").Value = 27500
$ws.Range("M65").Value = -24380
$ws.Range("H74").Value = 22140
$ws.Range("I74").Value = 21328.334
$ws.Range("J74").Value = 22748.75
$ws.Range("K74").Value = 21328.334
$ws.Range("L74").Value = 22748.75
$ws.Range("M74").Value = -20392.334
$ws.Range("N74").Value = -24620.75
$ws.Range("H77").Value = 22140
$ws.Range("I77").Value = 21328.334
$ws.Range("J77").Value = 22748.75
$ws.Range("K77").Value = 63985.00199999999
$ws.Range("L77").Value = 68246.25
$ws.Range("M77").Value = -59305.00199999999
$ws.Range("N77").Value = -77606.25
$ws.Range("H108").Value = 80625
$ws.Range("J108").Value = 80625
$ws.Range("L108").Value = 80625
$ws.Range("N108").Value = -88305
$ws.Range("H112").Value = 39193.5
$ws.Range("J112").Value = 39193.5
$ws.Range("L112").Value = 39193.5
$ws.Range("N112").Value = -42147.5
$ws.Range("H132").Value = 250052320
$ws.Range("I132").Value = 103129
$ws.Range("K132").Value = 309387
$ws.Range("M132").Value = -306857

Write-Output "Applied 196 cell updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets."